# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet (row 2 values updated)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 193
$wsOff.Range("C2").Value = 134
$wsOff.Range("D2").Value = 41
$wsOff.Range("E2").Value = 21

# DEF sheet (row 2 values updated)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 142
$wsDef.Range("C2").Value = 95
$wsDef.Range("D2").Value = 34
$wsDef.Range("E2").Value = 23
